$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, col B)
$meta.Cells.Item(3, 2).Value2 = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, col B)
$meta.Cells.Item(8, 2).Value2 = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" (old row 11)
$meta.Rows.Item(11).Insert()
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Cells.Item(11, 1).Value2 = "Jurisdiction"
$meta.Cells.Item(11, 2).Value2 = ""

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Add the II-1 invariant constraint text to Birthplace.typeId (row 5), column AJ (Constraint(s))
$elements.Cells.Item(5, 36).Value2 = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
